# Apply the "latest code" update described in the commit:
#  - User Information sheet gets refreshed sample data
#      (First/Last name changed, "Zip Code" renamed to "Postal Code",
#       and the postal code value is now entered as a text value "1939"
#       instead of a number, so Excel flags it with the quote-prefix style)
#  - "User Information" becomes the active/selected sheet instead of
#      "Login Details", with the selection parked on the new postal-code cell.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("User Information")

# NOTE: write the two name cells before the header/value cells so the
# shared-string table ends up ordered the same way Excel produced it
# (Zamantuli, Xulu, then Postal Code, then 1939).
$ws2.Range("A2").Value = "Zamantuli"
$ws2.Range("B2").Value = "Xulu"
$ws2.Range("C1").Value = "Postal Code"

# Enter the postal code as text (leading apostrophe = quote-prefixed text
# entry) rather than a number, matching the new quotePrefix cell style.
$ws2.Range("C2").Value = "'1939"

# Make "User Information" the active sheet / tab, with C2 selected.
$ws2.Activate() | Out-Null
$ws2.Range("C2").Select() | Out-Null

# Best-effort: tell Excel to ignore the "number stored as text" warning on
# C2 (maps to <ignoredErrors> in the saved file when supported).
try {
    $errs = $ws2.Range("C2").Errors
    $numberAsText = $errs.Item(6)
    $numberAsText.Ignore = $true
} catch {
}
